$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TwoxTwoOutTax_1-2")
$ws.Activate()

# Header row (row 1): new/changed shared-string labels for S1:V1
$ws.Range("S1").Value = "'O4,I=0.1"
$ws.Range("T1").Value = "'ITAX=0.1"
$ws.Range("U1").Value = "'O2=.3,I=.2"
$ws.Range("V1").Value = "'ITAX=100%"

# Row 2
$ws.Range("S2").Value = 0.63672356266481245
$ws.Range("T2").Value = 0.54435839324634638
$ws.Range("U2").Value = 0.26298581656113768

# Row 3
$ws.Range("T3").Value = 1.4510658620690864
$ws.Range("U3").Value = 1.7255770390585923
$ws.Range("V3").Value = 1.9797958977022081

# Row 4
$ws.Range("T4").Value = 0.96551294639121377
$ws.Range("U4").Value = 0.91527149061040103
$ws.Range("V4").Value = 0.85094165634213248

# Row 5
$ws.Range("T5").Value = 1.253716855956764
$ws.Range("U5").Value = 1.3799838239044699
$ws.Range("V5").Value = 1.5055868687041585

# Row 6
$ws.Range("T6").Value = 0.77493692216806154
$ws.Range("U6").Value = 0.68107531417393674
$ws.Range("V6").Value = 0.59749254327701107

# Row 7
$ws.Range("T7").Value = 1
$ws.Range("U7").Value = 1
$ws.Range("V7").Value = 1

# Row 8
$ws.Range("T8").Value = 0.9292784514896667
$ws.Range("U8").Value = 0.92558632011647268
$ws.Range("V8").Value = 0.93690421725376627

# Row 9
$ws.Range("T9").Value = 0.84027401688006276
$ws.Range("U9").Value = 0.79214910866172539
$ws.Range("V9").Value = 0.76497908901895417

# Row 10
$ws.Range("T10").Value = 193.10258927785549
$ws.Range("U10").Value = 183.05429809302706
$ws.Range("V10").Value = 170.18833062727202

# Row 11
$ws.Range("S11").Value = 81.862278102834352
$ws.Range("T11").Value = 83.433232535765256
$ws.Range("U11").Value = 87.427869776331136
$ws.Range("V11").Value = "Undf"

# Row 12
$ws.Range("S12").Value = 18.017820510637556
$ws.Range("T12").Value = 16.115952567509943
$ws.Range("U12").Value = 9.4388503851486849
$ws.Range("V12").Value = "Undf"

# Row 13
$ws.Range("T13").Value = 28.125819289987454
$ws.Range("U13").Value = 31.827759090372432
$ws.Range("V13").Value = 35.02889640415664

# Row 14
$ws.Range("T14").Value = 69.539579787829297
$ws.Range("U14").Value = 62.832913390522151
$ws.Range("V14").Value = 55.604906859264212

# Row 15
$ws.Range("S15").Value = 37.149389102107399
$ws.Range("T15").Value = 36.501668500052006
$ws.Range("U15").Value = 34.695155696294023
$ws.Range("V15").Value = "Undf"

# Row 16
$ws.Range("S16").Value = 62.968805282449473
$ws.Range("T16").Value = 63.678287300094347
$ws.Range("U16").Value = 65.731733209666629
$ws.Range("V16").Value = "Undf"

# Row 17
$ws.Range("T17").Value = 55.221484068689186
$ws.Range("U17").Value = 52.663928705525308
$ws.Range("V17").Value = 50.510257155949404

# Row 18
$ws.Range("T18").Value = 45.026343426821327
$ws.Range("U18").Value = 47.933812623576941
$ws.Range("V18").Value = 50.510257291375197

# Row 19
$ws.Range("T19").Value = 200
$ws.Range("U19").Value = 200
$ws.Range("V19").Value = 200.00000000000003

# Row 20
$ws.Range("T20").Value = 89.310036642418936
$ws.Range("U20").Value = 85.126151992605983
$ws.Range("V20").Value = 81.498026207158986

# Row 21
$ws.Range("T21").Value = 113.59698981563102
$ws.Range("U21").Value = 121.17204299879612
$ws.Range("V21").Value = 129.37005269973994

# Row 22
$ws.Range("T22").Value = 193.10258927785549
$ws.Range("U22").Value = 183.05429809302706
$ws.Range("V22").Value = 170.18833062727202

# Row 23
$ws.Range("T23").Value = 0.96551294638927743
$ws.Range("U23").Value = 0.91527149046513534
$ws.Range("V23").Value = 0.85094165313636017
